$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(7.2666666666666604, 0.68960536218590596, 6.5565029348161801, 1.1476770824495399),
    @(6.9428571428571404, 0.48360888228030702, 6.4186928104575101, 1.2155838279976701),
    @(7.0124999999999904, 0.45672064766112702, 6.4297417631344498, 1.2958790169746499),
    @(6.3388888888888797, 1.3338077859558799, 6.1940327237728496, 1.3940071801028899),
    @(6.5444444444444398, 1.1535238549551601, 6.2612708018154297, 1.3128876677111201),
    @(6.4666666666666597, 1.55313303429623, 6.5571428571428498, 1.5362291495737199),
    @(6.6099999999999897, 1.1707594875881899, 6.4970846306255599, 1.14754688229681)
)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]
    $ws.Cells.Item($row, 5).Value = $values[$i][3]
}

$ws.Range("C16").Select()
